# Atualização de bases das ligas, do dia: 22-05-2024 às 20:16
#
# The scraped source data for several fixtures came back in a different
# row order than before; this swaps the full data payload (columns B..AB,
# i.e. everything except the running index in column A) between the two
# affected rows for each mismatched pair, and applies a handful of
# standalone odds corrections for the still-unplayed fixtures at the
# bottom of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($row1, $row2, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# Columns B (2) through AB (28) -- column A (the sequential id) is untouched.
$firstCol = 2
$lastCol = 28

Swap-Rows 88 89 $firstCol $lastCol
Swap-Rows 114 115 $firstCol $lastCol
Swap-Rows 163 164 $firstCol $lastCol
Swap-Rows 195 196 $firstCol $lastCol
Swap-Rows 197 198 $firstCol $lastCol
Swap-Rows 208 209 $firstCol $lastCol
Swap-Rows 217 218 $firstCol $lastCol

# Standalone odds corrections (not part of a row swap) for the two
# not-yet-played fixtures at the end of the sheet.
$ws.Range("M311").Value = 2.05
$ws.Range("O311").Value = 3.5
$ws.Range("P311").Value = -0.5
$ws.Range("Q311").Value = 2.05
$ws.Range("R311").Value = 1.85
$ws.Range("S311").Value = 2.75
$ws.Range("T311").Value = 1.975
$ws.Range("U311").Value = 1.875

$ws.Range("Q312").Value = 2.02
$ws.Range("R312").Value = 1.88
$ws.Range("S312").Value = 2.5
$ws.Range("T312").Value = 1.8
$ws.Range("U312").Value = 2.05
